# This script re-orders the comma-separated "Recorded By" email lists in
# column G for a set of rows in the active worksheet, to match the
# canonical OOXML target. The underlying set of email addresses in each
# row is unchanged -- only the order within the comma-separated string
# is updated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = "majorelle.magdy@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, servinaz@med.asu.edu.eg"
    3  = "asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"
    4  = "asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg"
    5  = "eman.tantawi@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg, hananragab@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, servinaz@med.asu.edu.eg"
    6  = "asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, gehanadel@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, servinaz@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg"
    7  = "asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, gehanadel@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, servinaz@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg"
    8  = "majorelle.magdy@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, manar.montaser@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, backup@backdoor.com, shaimaa.ahmed@med.asu.edu.eg"
    9  = "majorelle.magdy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, manar.montaser@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, gehanadel@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg"
    10 = "alshimaa.atef@med.asu.edu.egm, gehanadel@med.asu.edu.eg, heba@med.asu.edu.eg, drsarahmw@gmail.com, shaimaa.ahmed@med.asu.edu.eg, servinaz@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg"
    11 = "asmaa.reda@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"
    12 = "salma.elgendy.std@med.asu.edu.eg, System"
    13 = "Omnia.Mohammed@med.asu.edu.eg, Shimaa.ashraf@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, wessam.atef@med.asu.edu.eg"
    17 = "marwa_mustafa@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, Sarah.Abdelmohsen@med.asu.edu.eg, dina.adel@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, nourhan.osama@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg"
    20 = "user@user.com, nourhan.mostafa@med.asu.edu.eg"
    22 = "Amr-Saeed@med.asu.edu.eg, nancy.abdelshafy@med.asu.edu.eg"
    23 = "hana.amr@med.asu.edu.eg, nourhan.mostafa@med.asu.edu.eg"
    24 = "marina_atef@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg"
    25 = "marina_atef@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, eman.samir@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg"
    26 = "youstina.magdy@med.asu.edu.eg, gigi@med.asu.edu.eg"
    27 = "eman.samir@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, Eman_mohamed@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg"
    28 = "abdullah.elagrody@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, eman.samir@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg"
    29 = "Monica.Eshak@med.asu.edu.eg, eman.samir@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg"
    30 = "asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg"
    31 = "asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"
    32 = "asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg"
    33 = "eman.tantawi@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg, hananragab@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, servinaz@med.asu.edu.eg"
    34 = "asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, gehanadel@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, servinaz@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg"
    35 = "asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, gehanadel@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, servinaz@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg"
    36 = "majorelle.magdy@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, manar.montaser@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, backup@backdoor.com, shaimaa.ahmed@med.asu.edu.eg"
    37 = "majorelle.magdy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, manar.montaser@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, gehanadel@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg"
    38 = "alshimaa.atef@med.asu.edu.egm, gehanadel@med.asu.edu.eg, heba@med.asu.edu.eg, drsarahmw@gmail.com, shaimaa.ahmed@med.asu.edu.eg, servinaz@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg"
    39 = "asmaa.reda@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"
    40 = "salma.elgendy.std@med.asu.edu.eg, System"
    41 = "Omnia.Mohammed@med.asu.edu.eg, Shimaa.ashraf@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, wessam.atef@med.asu.edu.eg"
    45 = "marwa_mustafa@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, Sarah.Abdelmohsen@med.asu.edu.eg, dina.adel@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, nourhan.osama@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg"
    48 = "user@user.com, nourhan.mostafa@med.asu.edu.eg"
    50 = "Amr-Saeed@med.asu.edu.eg, nancy.abdelshafy@med.asu.edu.eg"
    51 = "hana.amr@med.asu.edu.eg, nourhan.mostafa@med.asu.edu.eg"
    52 = "marina_atef@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg"
    53 = "marina_atef@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, eman.samir@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg"
    54 = "youstina.magdy@med.asu.edu.eg, gigi@med.asu.edu.eg"
    55 = "eman.samir@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, Eman_mohamed@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg"
    56 = "abdullah.elagrody@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, eman.samir@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg"
    57 = "Monica.Eshak@med.asu.edu.eg, eman.samir@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg"
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 7).Value = $values[$row]
}
